# Applies the update described in the commit: append additional SMS/status
# log entries from Freddy Velez to the "Sheet" worksheet (column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

$newValues = @(
    "20 13:09>>> 4FD889D140   Freddy Velez",
    "20 13:11>>> 4FD889D140   Freddy Velez",
    "20 13:11>>> 4FD889D140   Freddy Velez",
    "20 15:57>>> 4FD889D140   Freddy Velez",
    "20 20:56>>> 4FD889D140   Freddy Velez",
    "20 20:56>>> 4FD889D140   Freddy Velez",
    "20 20:57>>> 4FD889D140   Freddy Velez",
    "20 20:58>>> 4FD889D140   Freddy Velez",
    "20 21:43>>> 4FD889D140   Freddy Velez",
    "20 21:43>>> 4FD889D140   Freddy Velez",
    "20 21:43>>> 4FD889D140   Freddy Velez",
    "20 21:43>>> 4FD889D140   Freddy Velez",
    "20 22:50>>> 4FD889D140   Freddy Velez",
    "20 22:50>>> 4FD889D140   Freddy Velez",
    "20 22:50>>> 4FD889D140   Freddy Velez",
    "20 23:04>>> 4FD889D140   Freddy Velez"
)

$startRow = 311
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}
